# Enhance build script optimization
# Append a new logged row (row 40) to each of the four ROW*-LIFTER
# worksheets, mirroring the existing row layout:
#   time | total-length(hex) | ID(hex) | actual-length(hex) | checksum(hex)
#   | total-length(dec) | ID(dec) | actual-length(dec) | checksum(dec)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: ROW35-FE-LIFTER
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$t1 = [double]"45744.32396019676"
$id1 = [double]"5.68631262647114e+23"
$ws1.Cells.Item(40, 1).Value = $t1
$ws1.Cells.Item(40, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(40, 2).Value = "0x01,0x90"
$ws1.Cells.Item(40, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(40, 4).Value = "0x01,0x7a"
$ws1.Cells.Item(40, 5).Value = "0xd"
$ws1.Cells.Item(40, 6).Value = 400
$ws1.Cells.Item(40, 7).Value = $id1
$ws1.Cells.Item(40, 8).Value = 378
$ws1.Cells.Item(40, 9).Value = 13

# ---------------------------------------------------------------
# Sheet 2: ROW35-MID-LIFTER
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$t2 = [double]"45744.17079440972"
$id2 = [double]"5.68631262647114e+23"
$ws2.Cells.Item(40, 1).Value = $t2
$ws2.Cells.Item(40, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(40, 2).Value = "0x01,0x90"
$ws2.Cells.Item(40, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(40, 4).Value = "0x01,0x7a"
$ws2.Cells.Item(40, 5).Value = "0xe"
$ws2.Cells.Item(40, 6).Value = 400
$ws2.Cells.Item(40, 7).Value = $id2
$ws2.Cells.Item(40, 8).Value = 378
$ws2.Cells.Item(40, 9).Value = 14

# ---------------------------------------------------------------
# Sheet 3: ROW02-FE-LIFTER
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$t3 = [double]"45744.31755600694"
$id3 = [double]"5.68631262647114e+23"
$ws3.Cells.Item(40, 1).Value = $t3
$ws3.Cells.Item(40, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws3.Cells.Item(40, 2).Value = "0x01,0x90"
$ws3.Cells.Item(40, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(40, 4).Value = "0x01,0x7a"
$ws3.Cells.Item(40, 5).Value = "0x3"
$ws3.Cells.Item(40, 6).Value = 400
$ws3.Cells.Item(40, 7).Value = $id3
$ws3.Cells.Item(40, 8).Value = 378
$ws3.Cells.Item(40, 9).Value = 3

# ---------------------------------------------------------------
# Sheet 4: ROW02-MID-LIFTER
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$t4 = [double]"45744.37402591435"
$id4 = [double]"9.85046333984776e+23"
$ws4.Cells.Item(40, 1).Value = $t4
$ws4.Cells.Item(40, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws4.Cells.Item(40, 2).Value = "0x01,0x90"
$ws4.Cells.Item(40, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(40, 4).Value = "0x01,0x7a"
$ws4.Cells.Item(40, 5).Value = "0x3"
$ws4.Cells.Item(40, 6).Value = 400
$ws4.Cells.Item(40, 7).Value = $id4
$ws4.Cells.Item(40, 8).Value = 378
$ws4.Cells.Item(40, 9).Value = 3
